$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$currencyFmt = "[$£-809]#,##0.00;[RED]\-[$£-809]#,##0.00"

# --- K8 formula update (extend SUM range to B100) ---
$ws.Range("K8").Formula = "=SUM(B3:B100)"

# --- Fill in G column formulas for rows 33-44, matching the existing
#     money formatting used throughout the rest of column G ---
for ($r = 33; $r -le 44; $r++) {
    $cell = $ws.Range("G$r")
    $cell.Formula = "=B$r*`$K`$6"
    $cell.NumberFormat = $currencyFmt
    $cell.HorizontalAlignment = -4108
}

# --- Row 42: add Manhours (B), Description (E) values ---
$ws.Range("B42").Value = 0.5
$ws.Range("E42").Value = "Create warning/over identifiers for budget"

# --- Row 43: add Manhours (B), Description (E), Comment (F) values ---
$ws.Range("B43").Value = 2
$ws.Range("E43").Value = "Create Sort Function for Transaction class"
$ws.Range("F43").Value = "SQL was not happy with this at all – would not accept params no matter what I did"

# --- Row 44: brand new row of data ---
$ws.Range("A44").Value = $ws.Range("A43").Value2
$ws.Range("B44").Value = 0.25
$ws.Range("C44").Value = $ws.Range("C43").Value2
$ws.Range("D44").Value = 2
$ws.Range("E44").Value = "Create Sort Function for Tag Class"
$ws.Range("F44").Value = $ws.Range("F19").Value2

# --- View state: scroll position + active selection ---
$win = $excel.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 1
$ws.Range("K9").Select()
